$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("30:31").Insert()
$ws.Range("A29:Q29").Copy()
$ws.Range("A30:Q30").PasteSpecial(-4122)
$origL = $ws.Range("L30").NumberFormat()
$origP = $ws.Range("P30").NumberFormat()
Write-Host ("origL=" + $origL + " origP=" + $origP)
$ws.Range("L30").NumberFormat = "@"
$ws.Range("L30").Value = "0"
$ws.Range("L30").NumberFormat = $origL
$ws.Range("P30").NumberFormat = "@"
$ws.Range("P30").Value = "27.0000"
$ws.Range("P30").NumberFormat = $origP
Write-Host "done"
